$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(39.78, 0.16, 70.68860268592834, 138.890625),
    @(39.78, 0.16, 245.2299783229828, 149.359375),
    @(39.78, 0.16, 74.51772165298462, 128.265625),
    @(39.78, 0.16, 79.88404107093811, 137.375),
    @(39.78, 0.16, 73.88481378555298, 120.9375)
)

$startRow = 91
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
